$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting for the new rows by copying row 272 styles down to 273:277
$ws.Range("A272:R272").Copy($ws.Range("A273:R277"))

# Row 273
$ws.Range("A273").Value = 271
$ws.Range("B273").Value = 44515
$ws.Range("C273").Value = 1017.630004882812
$ws.Range("D273").Value = 1031.97998046875
$ws.Range("E273").Value = 978.5999755859375
$ws.Range("F273").Value = 1013.390014648438
$ws.Range("G273").Value = 1013.390014648438
$ws.Range("H273").Value = 34775600
$ws.Range("I273").Value = "TSLA"
$ws.Range("J273").Value = -20.030029296875
$ws.Range("K273").Value = 16.02571541922433
$ws.Range("L273").Value = 15.89999825613839
$ws.Range("M273").Value = 50.1968901374677
$ws.Range("N273").Value = 591954500
$ws.Range("O273").Value = -2702170
$ws.Range("P273").Value = 9173734.151229439
$ws.Range("Q273").Value = -21.902978515625
$ws.Range("R273").Value = 15.69155339325325

# Row 274
$ws.Range("A274").Value = 272
$ws.Range("B274").Value = 44516
$ws.Range("C274").Value = 1003.309997558594
$ws.Range("D274").Value = 1057.199951171875
$ws.Range("E274").Value = 1002.179992675781
$ws.Range("F274").Value = 1054.72998046875
$ws.Range("G274").Value = 1054.72998046875
$ws.Range("H274").Value = 26542400
$ws.Range("I274").Value = "TSLA"
$ws.Range("J274").Value = 41.3399658203125
$ws.Range("K274").Value = 16.32249886648995
$ws.Range("L274").Value = 15.89999825613839
$ws.Range("M274").Value = 50.6555988022028
$ws.Range("N274").Value = 565412100
$ws.Range("O274").Value = -769710
$ws.Range("P274").Value = 10013620.44794822
$ws.Range("Q274").Value = -5.474987792968751
$ws.Range("R274").Value = 8.280853346972963

# Row 275
$ws.Range("A275").Value = 273
$ws.Range("B275").Value = 44517
$ws.Range("C275").Value = 1063.510009765625
$ws.Range("D275").Value = 1119.640014648438
$ws.Range("E275").Value = 1055.5
$ws.Range("F275").Value = 1089.010009765625
$ws.Range("G275").Value = 1089.010009765625
$ws.Range("H275").Value = 31445400
$ws.Range("I275").Value = "TSLA"
$ws.Range("J275").Value = 34.280029296875
$ws.Range("K275").Value = 16.92999703543527
$ws.Range("L275").Value = 15.89999825613839
$ws.Range("M275").Value = 51.56868554221396
$ws.Range("N275").Value = 596857500
$ws.Range("O275").Value = 8199570
$ws.Range("P275").Value = 6328406.014690166
$ws.Range("Q275").Value = -7.65599365234375
$ws.Range("R275").Value = 7.06714795926534

# Row 276
$ws.Range("A276").Value = 274
$ws.Range("B276").Value = 44518
$ws.Range("C276").Value = 1106.550048828125
$ws.Range("D276").Value = 1112
$ws.Range("E276").Value = 1075.02001953125
$ws.Range("F276").Value = 1096.380004882812
$ws.Range("G276").Value = 1096.380004882812
$ws.Range("H276").Value = 20898900
$ws.Range("I276").Value = "TSLA"
$ws.Range("J276").Value = 7.3699951171875
$ws.Range("K276").Value = 15.69821166992188
$ws.Range("L276").Value = 15.89999825613839
$ws.Range("M276").Value = 49.68069933915767
$ws.Range("N276").Value = 575958600
$ws.Range("O276").Value = 13795500
$ws.Range("P276").Value = 5541185.075258421
$ws.Range("Q276").Value = 7.23099365234375
$ws.Range("R276").Value = 9.68794853425962

# Row 277
$ws.Range("A277").Value = 275
$ws.Range("B277").Value = 44519
$ws.Range("C277").Value = 1098.869995117188
$ws.Range("D277").Value = 1138.719970703125
$ws.Range("E277").Value = 1092.699951171875
$ws.Range("F277").Value = 1137.06005859375
$ws.Range("G277").Value = 1137.06005859375
$ws.Range("H277").Value = 21168000
$ws.Range("I277").Value = "TSLA"
$ws.Range("J277").Value = 40.6800537109375
$ws.Range("K277").Value = 12.71607317243304
$ws.Range("L277").Value = 15.89999825613839
$ws.Range("M277").Value = 44.43682356669267
$ws.Range("N277").Value = 597126600
$ws.Range("O277").Value = 4246240
$ws.Range("P277").Value = 5670909.244533308
$ws.Range("Q277").Value = 20.15399169921875
$ws.Range("R277").Value = 5.70334642338831
